$wb = $excel.ActiveWorkbook

# The "Swiss" sheet is used as the template for the new "Portugal" sheet:
# duplicate it, place the copy right after "Swiss", then adjust its content.
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy([System.Reflection.Missing]::Value, $swiss)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Portugal"

# The Swiss template's A1 note cell points at a stray shared-string entry
# ("S"); realign it with the note text used by the other sheets.
$newSheet.Range("A1").Value = "Note: Do not change the column/rows index "

# Update the market name / ticket reference cells.
# B4 is written first and B2 second so the two brand new shared-string
# entries land in the same order as in the target workbook
# (index N = "NGC-3479/T2407", index N+1 = "Portugal Market").
$newSheet.Range("B4").Value = "NGC-3479/T2407"
$newSheet.Range("B2").Value = "Portugal Market"

# Narrow the columns to match the Portugal layout.
$newSheet.Columns.Item(1).ColumnWidth = 23.666666666666668
$newSheet.Columns.Item(2).ColumnWidth = 16.333333333333332
$newSheet.Columns.Item(3).ColumnWidth = 14.666666666666666
$newSheet.Columns.Item(4).ColumnWidth = 14.166666666666666

# With the narrower columns the wrapped-text rows grow to two lines tall.
$newSheet.Rows.Item(3).RowHeight = 28.8
$newSheet.Rows.Item(4).RowHeight = 28.8
$newSheet.Rows.Item(5).RowHeight = 28.8

# Make the new sheet the active tab with the same selection state recorded
# in the target file.
$newSheet.Activate()
$newSheet.Range("B4:B5").Select()
